# Update the "Förändrad" (Changed) date column (C) for every data row
# on the active sheet from 2023-09-17 (45186) to 2023-09-19 (45188).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column C by scanning up from the bottom of the
# sheet (equivalent to pressing Ctrl+Up), which is reliable for a simple
# contiguous data table like this one (falls back to the known extent of
# 458 rows if something unexpected happens).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 458 }

$ws.Range("C2:C$lastRow").Value = 45188
